$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late" / data columns shift right
# by one: old N->O, O->P, P->Q), matching the "Variable Instalments" layout
# that adds an extra column to the repayment schedule sheet.
$ws.Columns("N").Insert()

# The inserted column picks up formatting from the left neighbour (column M)
# width of 11 characters, but without the "best fit" flag.
$ws.Columns("N").ColumnWidth = 11 - 5/6

# Make "Repayment schedule" the active/selected sheet and select K20 on it.
$ws.Activate()
$ws.Range("K20").Select()
